$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2 header cell: align its format with the rest of row 2 (font + centered) ---
# (done before the text edits below so the format copied is the original one)
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# --- Header / label text corrections ---
$ws.Range("B1").Value = "GW"
$ws.Range("B2").Value = "COM - Battery (Lead-acid) ELC Storage: DayNite (accompanying tech to represent power)"
$ws.Range("C2").Value = "COM - Battery (Li-ion) ELC Storage: DayNite (accompanying tech to represent power)"
$ws.Range("D2").Value = "RSD - Battery (Lead-acid) ELC Storage: DayNite (accompanying tech to represent power)"
$ws.Range("E2").Value = "RSD - Battery (Li-ion) ELC Storage: DayNite (accompanying tech to represent power)"

# --- Corrected BATS / H2 storage capacities (GW, previously mislabeled GWh) ---
$ws.Range("B3").Value = 3.887
$ws.Range("D3").Value = 3.8849999999999998

$ws.Range("B4").Value = 0.54500000000000004
$ws.Range("D4").Value = 0.54500000000000004

$ws.Range("B5").Value = 1.5740000000000001
$ws.Range("D5").Value = 1.575

$ws.Range("B6").Value = 0.095
$ws.Range("D6").Value = 0.216

$ws.Range("B7").Value = 1.048
$ws.Range("D7").Value = 1.335

$ws.Range("B8").Value = 18.471
$ws.Range("D8").Value = 19.434000000000001

$ws.Range("B9").Value = 3.0219999999999998
$ws.Range("D9").Value = 3.0219999999999998

$ws.Range("B10").Value = 0.47799999999999998
$ws.Range("D10").Value = 0.47799999999999998

$ws.Range("B11").Value = 0.52400000000000002
$ws.Range("D11").Value = 1.1319999999999999

$ws.Range("B12").Value = 0.42699999999999999
$ws.Range("D12").Value = 10.805999999999999

$ws.Range("B13").Value = 3.4550000000000001
$ws.Range("D13").Value = 3.4580000000000002

$ws.Range("D14").Value = 30.056000000000001

$ws.Range("B15").Value = 0.22600000000000001
$ws.Range("D15").Value = 0.24399999999999999

$ws.Range("B16").Value = 1.181
$ws.Range("D16").Value = 1.1839999999999999

$ws.Range("B17").Value = 1.478
$ws.Range("D17").Value = 1.478

$ws.Range("B18").Value = 0.96199999999999997
$ws.Range("D18").Value = 0.97199999999999998

$ws.Range("B19").Value = 5.8380000000000001
$ws.Range("D19").Value = 5.8380000000000001

$ws.Range("B20").Value = 0.47299999999999998
$ws.Range("D20").Value = 0.47299999999999998

$ws.Range("B21").Value = 0.042
$ws.Range("D21").Value = 0.042

$ws.Range("B22").Value = 9.6210000000000004
$ws.Range("D22").Value = 9.6210000000000004

$ws.Range("B23").Value = 1.675
$ws.Range("D23").Value = 1.8240000000000001

$ws.Range("B24").Value = 1.4339999999999999
$ws.Range("D24").Value = 2.605

$ws.Range("B25").Value = 1.994
$ws.Range("D25").Value = 3.7229999999999999

$ws.Range("D26").Value = 1.2070000000000001

$ws.Range("B27").Value = 5.86
$ws.Range("D27").Value = 5.86

$ws.Range("B28").Value = 0.55300000000000005
$ws.Range("D28").Value = 0.58399999999999996

$ws.Range("D29").Value = 0.216

$ws.Range("B30").Value = 14.048999999999999
$ws.Range("D30").Value = 14.048999999999999

# --- Number format: capacities now shown with two decimals (was custom #,##0.0) ---
$dataCells = "B3,D3,B4,D4,B5,C5,D5,E5,B6,D6,B7,D7,B8,D8,B9,D9,B10,D10,B11,D11,B12,D12,B13,D13,E13,D14,B15,D15,B16,D16,B17,C17,D17,E17,B18,D18,B19,D19,B20,D20,B21,D21,B22,D22,B23,D23,B24,D24,B25,D25,D26,B27,C27,D27,E27,B28,D28,D29,B30,D30"
$dataRange = $ws.Range($dataCells)
foreach ($area in $dataRange.Areas) {
    $area.NumberFormat = "#,##0.00"
}

# --- Widen columns B:E to fit the longer descriptive text ---
$ws.Range("B1").Select()
$ws.Columns("B").ColumnWidth = 67.1796875
$ws.Columns("C").ColumnWidth = 63.90625
$ws.Columns("D").ColumnWidth = 66.81640625
$ws.Columns("E").ColumnWidth = 63.54296875
